# The doc is a single-column table; each row is one stat value from the
# benchmark dump. Several values repeat elsewhere in the table, so we
# address cells by (row, column) instead of Find/Replace to avoid
# touching the wrong occurrence. Rows 44-46 originally packed a whole
# tab-separated line of per-run stats into one cell/run; the new content
# collapses each of those cells down to a single value, so rewriting the
# whole cell range (rather than editing individual <w:t> runs) is what
# naturally drops the old tabs/runs and leaves one clean <w:t>.
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text  = "0M"
$t.Cell(2,1).Range.Text  = "0M"
$t.Cell(3,1).Range.Text  = "0M"
$t.Cell(4,1).Range.Text  = "335"
$t.Cell(5,1).Range.Text  = "0.00003"
$t.Cell(6,1).Range.Text  = "0.00007"
$t.Cell(9,1).Range.Text  = "0.00005"
$t.Cell(10,1).Range.Text = "0.00005"
$t.Cell(12,1).Range.Text = "0.01436"
$t.Cell(44,1).Range.Text = "100"
$t.Cell(45,1).Range.Text = "0.01"
$t.Cell(46,1).Range.Text = "342"
